$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly record is inserted at row 111, pushing the existing rows
# 111-211 down to 112-212 (dimension grows from A1:R211 to A1:R212).
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row with the latest observation. All
# columns besides the date (D) and volume (J) repeat the values that
# used to sit in row 111 before the shift.
$ws.Cells.Item(111, 1).Value = 8
$ws.Cells.Item(111, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(111, 3).Value = "Coquimbo"
$ws.Cells.Item(111, 4).Value = 44589
$ws.Cells.Item(111, 5).Value = 4
$ws.Cells.Item(111, 6).Value = 100112003
$ws.Cells.Item(111, 7).Value = "Ajo"
$ws.Cells.Item(111, 8).Value = "Chino"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 500
$ws.Cells.Item(111, 11).Value = 19000
$ws.Cells.Item(111, 12).Value = 20000
$ws.Cells.Item(111, 13).Value = 19500
$ws.Cells.Item(111, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(111, 15).Value = "China"
$ws.Cells.Item(111, 16).Value = 1950
$ws.Cells.Item(111, 17).Value = 10
$ws.Cells.Item(111, 18).Value = "Hortaliza"
